$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "65.143.47"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "3.181.92"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.43%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.119"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.71%  "

$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").Value = "3.730.81"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").Value = "65.030.84"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.63"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "3.171.95"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000157"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "415.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.62%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.67"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("E24").Value = "  -1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.488"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000105"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.46"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.35"
$ws.Range("D32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.94"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("E35").Value = "  -1.99%  "

$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.73"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "2.724.75"
$ws.Range("E37").Value = "  -3.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.90"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.12"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.709"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0634"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "293.86"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0988"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.98"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.48"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.78"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.901"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.57%  "

